# Weekly price update: insert a new price record for Cilantro at
# "Feria Lagunitas de Puerto Montt" dated 2023-10-24 (serial 45223),
# pushing the existing rows 500-542 down to 501-543.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 500; everything below (old rows 500-542)
# shifts down to 501-543, and the sheet dimension grows to A1:R543.
$ws.Rows.Item(500).Insert()

# Populate the newly inserted row 500 with the new weekly record.
$ws.Range("A500").Value = 4
$ws.Range("B500").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C500").Value = "Los Lagos"
$ws.Range("D500").Value = 45223
$ws.Range("E500").Value = 10
$ws.Range("F500").Value = 100112040
$ws.Range("G500").Value = "Cilantro"
$ws.Range("H500").Value = "Sin especificar"
$ws.Range("I500").Value = "Primera"
$ws.Range("J500").Value = 180
$ws.Range("K500").Value = 14000
$ws.Range("L500").Value = 15000
$ws.Range("M500").Value = 14500
$ws.Range("N500").Value = "$/caja 36 atados"
$ws.Range("O500").Value = "Región Metropolitana"
$ws.Range("P500").Value = 403
$ws.Range("Q500").Value = 36
$ws.Range("R500").Value = "Hortaliza"
